# Apply updated coin price/volume data (symbol-list refresh, Jan 18 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.17%"

# Row 3: OKB
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.41%"

# Row 4: HuobiToken
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.956"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.74%"

# Row 5: Cronos
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07638"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.15%"

# Row 6: FTXToken
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.930"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-14.61%"

# Row 7: KuCoinToken
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.833"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.42%"

# Row 8: MXToken
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9174"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.24%"

# Row 9: WazirX
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1748"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.12%"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07751"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.00%"

# Row 11: MandalaExchangeToken
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08530"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.94%"

# Row 12: BitrueCoin
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03239"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.58%"

# Row 13: BitMartToken
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.05%"

# Row 14: BitForexToken
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001511"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"

# Row 15: TigerCash
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005940"

# Row 16: LEO
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.464"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.40%"

# Row 17: GateToken
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.801"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.80%"

# Row 18: BTSEToken
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.37%"

# Row 19: BitpandaEcosystemToken
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3351"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.85%"

# Row 20: ProBitToken
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.84%"

# Row 21: MCDex
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.33%"

# Row 22: ZBToken
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.52%"

# Row 23: CoinExToken
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04524"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.63%"

# Row 24: BitKan
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.27%"

# Row 25: HotbitToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004388"

# Row 26: NitroEx
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.18%"

# Row 39: One
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.68%"

# Row 40: IDEX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04691"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.41%"

# Row 41: KickToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007496"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.95%"

# Row 42: BKEXToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.54%"

# Row 43: CEJI
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002334"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.58%"

# Row 44: LocalTraders
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01054"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.24%"

# Row 45: CoinLion
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006244"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.55%"

# Row 46: Kangarootoken
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.14%"

# Row 47: BOLO
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8204"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-28.70%"

# Row 49: CryptobidCoin
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.14%"

# Row 50: SpecialPowerGold
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
